# This script updates the "want to go" (F) and "minimum ticket price" (G)
# figures for a set of conventions that appear on both the "展览" (Exhibition)
# sheet and the "全部类型" (All types) sheet. In "全部类型" every affected
# row is shifted down by exactly one row relative to "展览" because that
# sheet contains one extra (演出) entry above them.

$wb = $excel.ActiveWorkbook

# Row (relative to the "展览" sheet) -> hashtable of column letter -> new value
$updates = @{
    3  = @{ F = 5579; G = 70 }
    5  = @{ F = 64 }
    6  = @{ G = 65 }
    7  = @{ F = 668 }
    8  = @{ F = 652 }
    9  = @{ F = 15 }
    10 = @{ F = 4 }
    11 = @{ F = 1084 }
    13 = @{ F = 1562 }
    14 = @{ F = 5157 }
    15 = @{ F = 458 }
    16 = @{ F = 263 }
    17 = @{ F = 228; G = 65 }
    18 = @{ F = 39 }
    19 = @{ F = 14 }
    21 = @{ F = 4459 }
    22 = @{ F = 223 }
    23 = @{ F = 1176 }
    25 = @{ F = 73 }
    27 = @{ F = 61 }
    28 = @{ F = 183 }
    29 = @{ F = 67 }
    32 = @{ F = 350 }
    33 = @{ F = 17 }
    34 = @{ F = 43 }
    35 = @{ F = 69 }
    36 = @{ F = 12 }
    37 = @{ F = 31 }
    40 = @{ F = 45 }
}

function Apply-Updates($ws, $rowOffset) {
    foreach ($row in $updates.Keys) {
        $targetRow = $row + $rowOffset
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$targetRow").Value = $cols[$col]
        }
    }
}

$wsExhibition = $wb.Worksheets.Item("展览")
Apply-Updates $wsExhibition 0

$wsAllTypes = $wb.Worksheets.Item("全部类型")
Apply-Updates $wsAllTypes 1
